# Apply "repull data, push all data, mean calculation" edit:
# Update column F (dSF) values for a set of rows on Sheet1 to reflect
# the re-pulled data values described in the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    14 = -1
    19 = -1
    23 = -1
    25 = 3
    37 = -4
    38 = 0
    44 = 4
    49 = 5
    55 = 0
    56 = -1
    60 = 1
    63 = -1
    66 = 6
    68 = 3
    69 = 1
    72 = 4
    76 = 4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
